$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week date range) ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Data table updates (rows 15-30) ---

# Row 15
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 3
$ws.Range("M15").Value = -100
$ws.Range("L15").Copy()
$ws.Range("M15").PasteSpecial(-4122)

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -20
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -88.235294117647

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 133.333333333333
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 53.846153846153
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 5
$ws.Range("K17").Value = 140
$ws.Range("L17").Value = 300
$ws.Range("M17").Value = 500
$ws.Range("N17").Value = 50

# Row 18
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 4
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -87.5
$ws.Range("N18").Value = -97.014925373134

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 116.666666666667
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -2.127659574468
$ws.Range("I19").Value = 24
$ws.Range("J19").Value = 22
$ws.Range("K19").Value = 9.090909090909
$ws.Range("L19").Value = -42.857142857142
$ws.Range("M19").Value = 26.315789473684
$ws.Range("N19").Value = 9.090909090909

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("C20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = 50
$ws.Range("H20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 2
$ws.Range("I20").Copy()
$ws.Range("J20").PasteSpecial(-4122)
$ws.Range("K20").Value = 250
$ws.Range("L20").Copy()
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("L20").Value = 40
$ws.Range("M20").Value = -12.5
$ws.Range("N20").Value = -88.888888888888

# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 47.058823529411
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = 18.823529411764
$ws.Range("I21").Value = 49
$ws.Range("J21").Value = 40
$ws.Range("K21").Value = 22.5
$ws.Range("L21").Value = -18.333333333333
$ws.Range("M21").Value = -5.769230769230
$ws.Range("N21").Value = -74.742268041237

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1

# Row 24
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 80
$ws.Range("F24").Value = 128
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = 19.626168224299
$ws.Range("I24").Value = 72
$ws.Range("J24").Value = 48
$ws.Range("K24").Value = 50
$ws.Range("L24").Value = 20
$ws.Range("M24").Value = 89.473684210526

# Row 25
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 96.551724137931
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = 112.5
$ws.Range("L25").Value = 209.090909090909
$ws.Range("M25").Value = 78.947368421052

# Row 26
$ws.Range("D26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("J26").Value = 3

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -50
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 2
$ws.Range("G27").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("K27").Value = 150
$ws.Range("H27").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("L27").Value = 400
$ws.Range("H27").Copy()
$ws.Range("L27").PasteSpecial(-4122)

# Row 30
$ws.Range("L30").Value = -100
$ws.Range("L16").Copy()
$ws.Range("L30").PasteSpecial(-4122)
